# ---------------------------------------------------------------------------
# Commit: "#5: fund, bonds, otherbonds, antique done"
#
# 1. Delete the "其他有價證券" (other securities) sheet entirely.
# 2. Rebuild the "基金受益憑證" (fund) sheet with the fuller record layout
#    (dealer column + the property_category/category/date/legislator_*/
#    source_file/index metadata columns used by the other "value" sheets),
#    matching the data already present in the workbook (quantities, face
#    values, totals, dealer banks, owners, currencies).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. remove the "其他有價證券" sheet -------------------------------------
$other = $wb.Worksheets.Item("其他有價證券")
$other.Delete()

# --- 2. rebuild "基金受益憑證" ----------------------------------------------
$ws  = $wb.Worksheets.Item("基金受益憑證")
$ref = $wb.Worksheets.Item("股票")   # donor sheet for styles / text-typed date

$ws.Cells.Clear()

# Header row (row 1): name / owner / dealer / quantity / face_value /
# currency / total / property_category / category / date / legislator_name /
# legislator_id / source_file / index
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Data rows 2-10
$ws.Range("A2").Value = 122
$ws.Range("B2").Value = "施羅德新興亞洲"
$ws.Range("C2").Value = "吳〇庭"
$ws.Range("D2").Value = "台北富邦商業銀行"
$ws.Range("E2").Value = 6036
$ws.Range("F2").Value = 21.62
$ws.Range("G2").Value = "美金"
$ws.Range("H2").Value = 39279.99
$ws.Range("I2").Value = "fund"
$ws.Range("J2").Value = "normal"
$ws.Range("L2").Value = "吳育昇"
$ws.Range("M2").Value = 1322
$ws.Range("N2").Value = "tmpe6fb1"
$ws.Range("O2").Value = 122

$ws.Range("A3").Value = 123
$ws.Range("B3").Value = "貝萊德新興歐洲"
$ws.Range("C3").Value = "吳〇庭"
$ws.Range("D3").Value = "台北富邦商業銀行"
$ws.Range("E3").Value = 8.06
$ws.Range("F3").Value = 110.82
$ws.Range("G3").Value = "美金"
$ws.Range("H3").Value = 26885.6
$ws.Range("I3").Value = "fund"
$ws.Range("J3").Value = "normal"
$ws.Range("L3").Value = "吳育昇"
$ws.Range("M3").Value = 1322
$ws.Range("N3").Value = "tmpe6fb1"
$ws.Range("O3").Value = 123

$ws.Range("A4").Value = 124
$ws.Range("B4").Value = "富達歐洲"
$ws.Range("C4").Value = "吳〇學"
$ws.Range("D4").Value = "台北富邦商業銀行"
$ws.Range("E4").Value = 85.83
$ws.Range("F4").Value = 8.44
$ws.Range("G4").Value = "歐元"
$ws.Range("H4").Value = 28744.4
$ws.Range("I4").Value = "fund"
$ws.Range("J4").Value = "normal"
$ws.Range("L4").Value = "吳育昇"
$ws.Range("M4").Value = 1322
$ws.Range("N4").Value = "tmpe6fb1"
$ws.Range("O4").Value = 124

$ws.Range("A5").Value = 125
$ws.Range("B5").Value = "貝萊德拉丁美洲"
$ws.Range("C5").Value = "劉娟娟"
$ws.Range("D5").Value = "國泰世華商業銀行"
$ws.Range("E5").Value = 27.66
$ws.Range("F5").Value = 81.09
$ws.Range("G5").Value = "美金"
$ws.Range("H5").Value = 67512.78
$ws.Range("I5").Value = "fund"
$ws.Range("J5").Value = "normal"
$ws.Range("L5").Value = "吳育昇"
$ws.Range("M5").Value = 1322
$ws.Range("N5").Value = "tmpe6fb1"
$ws.Range("O5").Value = 125

$ws.Range("A6").Value = 126
$ws.Range("B6").Value = "摩根東協"
$ws.Range("C6").Value = "劉娟娟"
$ws.Range("D6").Value = "國泰世華商業銀行"
$ws.Range("E6").Value = 24.88
$ws.Range("F6").Value = 92.89
$ws.Range("G6").Value = "美金"
$ws.Range("H6").Value = 69564.21
$ws.Range("I6").Value = "fund"
$ws.Range("J6").Value = "normal"
$ws.Range("L6").Value = "吳育昇"
$ws.Range("M6").Value = 1322
$ws.Range("N6").Value = "tmpe6fb1"
$ws.Range("O6").Value = 126

$ws.Range("A7").Value = 127
$ws.Range("B7").Value = "摩根大中華"
$ws.Range("C7").Value = "劉娟娟"
$ws.Range("D7").Value = "國泰世華商業銀行"
$ws.Range("E7").Value = 1208.587
$ws.Range("F7").Value = 23.37
$ws.Range("G7").Value = "美金"
$ws.Range("H7").Value = 850164.81
$ws.Range("I7").Value = "fund"
$ws.Range("J7").Value = "normal"
$ws.Range("L7").Value = "吳育昇"
$ws.Range("M7").Value = 1322
$ws.Range("N7").Value = "tmpe6fb1"
$ws.Range("O7").Value = 127

$ws.Range("A8").Value = 129
$ws.Range("B8").Value = "貝萊德新興歐洲"
$ws.Range("C8").Value = "劉娟娟"
$ws.Range("D8").Value = "台新國際商業銀行"
$ws.Range("E8").Value = 17.9
$ws.Range("F8").Value = 80.67
$ws.Range("G8").Value = "歐元"
$ws.Range("H8").Value = 57297.64
$ws.Range("I8").Value = "fund"
$ws.Range("J8").Value = "normal"
$ws.Range("L8").Value = "吳育昇"
$ws.Range("M8").Value = 1322
$ws.Range("N8").Value = "tmpe6fb1"
$ws.Range("O8").Value = 129

$ws.Range("A9").Value = 130
$ws.Range("B9").Value = "坦伯頓全球亞洲成長"
$ws.Range("C9").Value = "劉娟娟"
$ws.Range("D9").Value = "台新國際商業銀行"
$ws.Range("E9").Value = 74.308
$ws.Range("F9").Value = 28.67
$ws.Range("G9").Value = "美金"
$ws.Range("H9").Value = 64125.35
$ws.Range("I9").Value = "fund"
$ws.Range("J9").Value = "normal"
$ws.Range("L9").Value = "吳育昇"
$ws.Range("M9").Value = 1322
$ws.Range("N9").Value = "tmpe6fb1"
$ws.Range("O9").Value = 130

$ws.Range("A10").Value = 131
$ws.Range("B10").Value = "富達拉丁美洲"
$ws.Range("C10").Value = "劉娟娟"
$ws.Range("D10").Value = "台新國際商業銀行"
$ws.Range("E10").Value = 52.17
$ws.Range("F10").Value = 40.17
$ws.Range("G10").Value = "美金"
$ws.Range("H10").Value = 63079.63
$ws.Range("I10").Value = "fund"
$ws.Range("J10").Value = "normal"
$ws.Range("L10").Value = "吳育昇"
$ws.Range("M10").Value = 1322
$ws.Range("N10").Value = "tmpe6fb1"
$ws.Range("O10").Value = 131

# Column K ("date") holds the literal text "2011-11-17" on every data row,
# same as every other sheet in this workbook. Assigning that string directly
# via .Value gets auto-coerced to a date serial, so instead copy the already
# text-typed date column from the "股票" (stock) sheet, which carries the
# exact same value on every one of its 9 data rows.
$ref.Range("J2:J10").Copy()
$ws.Range("K2:K10").PasteSpecial(-4163)

# Match the workbook's existing look: header row + the "index" column (A)
# use the bold/bordered style that every other sheet's header & index column
# already uses - copy it over instead of re-deriving new style entries.
$ref.Range("B1").Copy()
$ws.Range("B1:O1").PasteSpecial(-4122)
$ref.Range("A2").Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)
